# Hypatia teaching workbook - add "Max_land_usage" parameter sheet
# (mirrors the existing "Max_newcap" sheet layout, but with 1E+20 values
# and a distinct bold style for the header/label cells).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new worksheet as the very last tab (after the current
#    last sheet, "Min_RES_electric_penetration") and name it.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Max_land_usage"

# ------------------------------------------------------------------
# 2. Header row (row 1): grouped headers with merged cells.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Tech_category"
$ws.Range("B1").Value = "Supply"
$ws.Range("G1").Value = "Conversion"
$ws.Range("I1").Value = "Transmission"

$ws.Range("B1:F1").Merge()
$ws.Range("G1:H1").Merge()

# ------------------------------------------------------------------
# 3. Row 2: technology names.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Technology"
$ws.Range("B2").Value = "Natural_gas_supply"
$ws.Range("C2").Value = "Oil_supply"
$ws.Range("D2").Value = "PV_PP"
$ws.Range("E2").Value = "Wind_PP"
$ws.Range("F2").Value = "Hydro_PP"
$ws.Range("G2").Value = "HFO_PP"
$ws.Range("H2").Value = "OCGT_PP"
$ws.Range("I2").Value = "Elec_transmission_distribution"

# ------------------------------------------------------------------
# 4. Row 3: "Years" label.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "Years"

# ------------------------------------------------------------------
# 5. Rows 4-14: Y0..Y10 labels in column A, 1E+20 across B..I.
# ------------------------------------------------------------------
$years = @("Y0","Y1","Y2","Y3","Y4","Y5","Y6","Y7","Y8","Y9","Y10")
$bigVal = [double]"1E+20"
for ($i = 0; $i -lt $years.Length; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $years[$i]
    for ($c = 2; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $bigVal
    }
}

# ------------------------------------------------------------------
# 6. Styling: bold font, thin box border, center/top alignment.
#    Two distinct (but visually identical) direct formats are used,
#    matching the source file:
#      - "label" style  -> column A (rows 1-14), I1, and all of row 2
#      - "merged" style -> the blank merged-header placeholders B1:H1
#    Applied as separate single-area ranges since a multi-area Range
#    object only formats its first area.
# ------------------------------------------------------------------
function Set-LabelStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Family = 0
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

Set-LabelStyle $ws.Range("A1:A14")
Set-LabelStyle $ws.Range("I1")
Set-LabelStyle $ws.Range("B2:I2")

function Set-MergedHeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Family = 0
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

Set-MergedHeaderStyle $ws.Range("B1:H1")

# ------------------------------------------------------------------
# 7. Page margins, matching the other parameter sheets in this
#    workbook (0.75in left/right, 1in top/bottom, 0.5in header/footer).
# ------------------------------------------------------------------
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 8. Sheet view bits: new sheet becomes the active/selected tab with
#    default selection, matching the source workbook state.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select()

Write-Host "Max_land_usage sheet created"
